$p = $ppt.ActivePresentation

# Insert a new slide "Der Kunde" right after slide 3 ("Die Lösung") and
# before the former slide 4 ("KANSROOSA"), using the "Titel und Inhalt"
# (Title and Content) layout - same layout used by the other content
# slides in this deck.
$newSlide = $p.Slides.Add(4, 2)

# Title placeholder
$title = $newSlide.Shapes.Item(1).TextFrame.TextRange
$title.LanguageID = "de-CH"
$title.Text = "Der Kunde"

# Body / content placeholder
$body = $newSlide.Shapes.Item(2).TextFrame.TextRange
$body.LanguageID = "de-CH"
$body.Text = "Fitnesscenter als primärer Kunde`r`rMitglieder profitieren vom Angebot`rFitnesscenter profitiert von zufrieden Kunden`r`r"
